$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'71.269.77"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +3.65%  "

$ws.Range("D3").Value = "'2.626.06"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +4.07%  "

$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("D5").Value = "'605.81"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.91%  "

$ws.Range("D6").Value = "'179.98"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.52%  "

$ws.Range("D8").Value = "'0.525"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.97%  "

$ws.Range("D9").Value = "'2.624.84"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.05%  "

$ws.Range("D10").Value = "'0.165"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +13.09%  "

$ws.Range("E11").Value = "  +0.25%  "

$ws.Range("D12").Value = "'0.347"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.62%  "

$ws.Range("D13").Value = "'5.00"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.42%  "

$ws.Range("D14").Value = "'3.111.36"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.38%  "

$ws.Range("D15").Value = "'0.0000186"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +9.14%  "

$ws.Range("D16").Value = "'26.72"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.15%  "

$ws.Range("D17").Value = "'71.267.07"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.82%  "

$ws.Range("D18").Value = "'2.624.01"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +4.30%  "

$ws.Range("D19").Value = "'381.43"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +8.11%  "

$ws.Range("D20").Value = "'7.87"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.97%  "

$ws.Range("D21").Value = "'11.47"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.39%  "

$ws.Range("D22").Value = "'4.12"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.48%  "

$ws.Range("E23").Value = "  +15.97%  "

$ws.Range("D24").Value = "'72.29"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.33%  "

$ws.Range("D25").Value = "'4.43"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.88%  "

$ws.Range("E26").Value = "  +0.18%  "

$ws.Range("D27").Value = "'9.92"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +9.80%  "

$ws.Range("D28").Value = "'2.763.46"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.03%  "

$ws.Range("D29").Value = "'1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.17%  "

$ws.Range("B30").Value = "Bittensor"
$ws.Range("C30").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D30").Value = "'548.37"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.25%  "

$ws.Range("B31").Value = "PEPE"
$ws.Range("C31").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D31").Value = "'0.0₃0958"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +7.68%  "

$ws.Range("D32").Value = "'8.07"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.50%  "

$ws.Range("E33").Value = "  +5.84%  "

$ws.Range("D34").Value = "'1.83"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.34%  "

$ws.Range("D35").Value = "'0.999"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.06%  "

$ws.Range("D36").Value = "'166.06"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.90%  "

$ws.Range("D38").Value = "'19.24"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.35%  "

$ws.Range("D39").Value = "'1.89"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +6.44%  "

$ws.Range("D40").Value = "'19.04"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.93%  "

$ws.Range("E41").Value = "  +5.47%  "

$ws.Range("B42").Value = "USDe"
$ws.Range("C42").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D42").Value = "'1.00"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.07%  "

$ws.Range("B43").Value = "dogwifhat"
$ws.Range("C43").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D43").Value = "'2.62"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +8.75%  "

$ws.Range("D44").Value = "'5.04"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.97%  "

$ws.Range("D45").Value = "'0.331"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.82%  "

$ws.Range("D46").Value = "'39.93"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.48%  "

$ws.Range("D47").Value = "'153.64"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.49%  "

$ws.Range("D48").Value = "'3.64"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.35%  "

$ws.Range("D49").Value = "'0.536"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.93%  "

$ws.Range("D50").Value = "'1.68"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.76%  "

$ws.Range("E51").Value = "  +2.00%  "

